# Applies the WR_89700562_WeekEnding_080325 update:
#  - refresh report-generated timestamp, total billed amount, line-item
#    count, and billing period
#  - repurpose the old row-66 TOTAL row into a data row and append the
#    rest of Tuesday/Point-04 + Point-03 items, closing with a new TOTAL
#  - append an entire new "Thursday (07/31/2025)" section (header, column
#    headers, two data rows, TOTAL)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header / summary cell updates
# ---------------------------------------------------------------------
$ws.Range("D5").Value2  = "Report Generated On: 08/16/2025 12:48 AM"
$ws.Range("C8").Value2  = 17736.03
$ws.Range("C9").Value2  = 47
$ws.Range("C10").Value2 = "07/28/2025 to 08/03/25"

# ---------------------------------------------------------------------
# 2. Build new rows 66-71 (tail end of the Tuesday/Point 04+03 table)
#    using the existing styled rows 64 (shaded) / 65 (plain) / 58 (TOTAL)
#    as formatting templates, then overwrite the cell text/values.
# ---------------------------------------------------------------------
$ws.Range("A64:H64").Copy($ws.Range("A66:H66"))
$ws.Range("A65:H65").Copy($ws.Range("A67:H67"))
$ws.Range("A64:H64").Copy($ws.Range("A68:H68"))
$ws.Range("A65:H65").Copy($ws.Range("A69:H69"))
$ws.Range("A64:H64").Copy($ws.Range("A70:H70"))
$ws.Range("A58").Copy($ws.Range("A71"))
$ws.Range("H58").Copy($ws.Range("H71"))

# Row 66
$ws.Range("A66").Value2 = "Point 04"
$ws.Range("B66").Value2 = "PLA-DLOC"
$ws.Range("C66").Value2 = "Inst"
$ws.Range("D66").Value2 = "PLA,Difficult Location"
$ws.Range("E66").Value2 = "EA"
$ws.Range("F66").Value2 = 2
$ws.Range("G66").Value2 = ""
$ws.Range("H66").Value2 = 238.2

# Row 67
$ws.Range("A67").Value2 = "Point 04"
$ws.Range("B67").Value2 = "CNA-TM"
$ws.Range("C67").Value2 = "Rem"
$ws.Range("D67").Value2 = "CNA,Trim Trees Minor"
$ws.Range("E67").Value2 = "EA"
$ws.Range("F67").Value2 = 5
$ws.Range("G67").Value2 = ""
$ws.Range("H67").Value2 = 539.6

# Row 68
$ws.Range("A68").Value2 = "Point 04"
$ws.Range("B68").Value2 = "PLA-BACK"
$ws.Range("C68").Value2 = "Inst"
$ws.Range("D68").Value2 = "Difficult Location Equip Adder-Backyard"
$ws.Range("E68").Value2 = "EA"
$ws.Range("F68").Value2 = 6
$ws.Range("G68").Value2 = ""
$ws.Range("H68").Value2 = 714.6

# Row 69
$ws.Range("A69").Value2 = "Point 03"
$ws.Range("B69").Value2 = "PLA-DLOC"
$ws.Range("C69").Value2 = "Inst"
$ws.Range("D69").Value2 = "PLA,Difficult Location"
$ws.Range("E69").Value2 = "EA"
$ws.Range("F69").Value2 = 2
$ws.Range("G69").Value2 = ""
$ws.Range("H69").Value2 = 238.2

# Row 70
$ws.Range("A70").Value2 = "Point 03"
$ws.Range("B70").Value2 = "CNA-TM"
$ws.Range("C70").Value2 = "Rem"
$ws.Range("D70").Value2 = "CNA,Trim Trees Minor"
$ws.Range("E70").Value2 = "EA"
$ws.Range("F70").Value2 = 3
$ws.Range("G70").Value2 = ""
$ws.Range("H70").Value2 = 323.76

# Row 71 - TOTAL for the Tuesday/Point 04+03 table.
# The old A66:G66 merge is superseded; remove it, merge the fresh
# A71:G71 range *before* stamping style (merging a still-blank range
# avoids the engine fanning the anchor cell's style out across the
# whole merged area), then paste in style (PasteSpecial formats only)
# and finally the literal values so only A71/H71 end up as real cells -
# matching the source report's minimal TOTAL-row cell layout.
$ws.Range("A66:G66").UnMerge()
$ws.Range("A71:G71").Merge()
$ws.Range("A58").Copy()
$ws.Range("A71").PasteSpecial(-4122)
$ws.Range("H58").Copy()
$ws.Range("H71").PasteSpecial(-4122)
$ws.Range("A71").Value2 = "TOTAL"
$ws.Range("H71").Value2 = 3807.32

# ---------------------------------------------------------------------
# 3. Build the new "Thursday (07/31/2025)" section: rows 74-78.
#    Templates: row 61 (day header), row 62 (column headers),
#    row 65 (plain data row), row 64 (shaded data row), row 58 (TOTAL).
# ---------------------------------------------------------------------

# Row 74 - day header (merge-then-paste-formats, see note above, so the
# row ends up with only the single A74 cell, same as row 39/61).
$ws.Range("A74:H74").Merge()
$ws.Range("A61").Copy()
$ws.Range("A74").PasteSpecial(-4122)
$ws.Range("A74").Value2 = "Thursday (07/31/2025)"

$ws.Range("A62:H62").Copy($ws.Range("A75:H75"))
$ws.Range("A65:H65").Copy($ws.Range("A76:H76"))
$ws.Range("A64:H64").Copy($ws.Range("A77:H77"))

# Row 75 - column headers (values already correct from template, but set
# explicitly to guarantee parity)
$ws.Range("A75").Value2 = "Point Number"
$ws.Range("B75").Value2 = "Billable Unit Code"
$ws.Range("C75").Value2 = "Work Type"
$ws.Range("D75").Value2 = "Unit Description"
$ws.Range("E75").Value2 = "Unit of Measure"
$ws.Range("F75").Value2 = "# Units"
$ws.Range("G75").Value2 = "N/A"
$ws.Range("H75").Value2 = "Pricing"

# Row 76
$ws.Range("A76").Value2 = "Point 03"
$ws.Range("B76").Value2 = "CNA-TM"
$ws.Range("C76").Value2 = "Rem"
$ws.Range("D76").Value2 = "CNA,Trim Trees Minor"
$ws.Range("E76").Value2 = "EA"
$ws.Range("F76").Value2 = 3
$ws.Range("G76").Value2 = ""
$ws.Range("H76").Value2 = 323.76

# Row 77
$ws.Range("A77").Value2 = "Point 03"
$ws.Range("B77").Value2 = "PLA-DLOC"
$ws.Range("C77").Value2 = "Inst"
$ws.Range("D77").Value2 = "PLA,Difficult Location"
$ws.Range("E77").Value2 = "EA"
$ws.Range("F77").Value2 = 2
$ws.Range("G77").Value2 = ""
$ws.Range("H77").Value2 = 238.2

# Row 78 - TOTAL for Thursday (same merge-first technique as row 71)
$ws.Range("A78:G78").Merge()
$ws.Range("A58").Copy()
$ws.Range("A78").PasteSpecial(-4122)
$ws.Range("H58").Copy()
$ws.Range("H78").PasteSpecial(-4122)
$ws.Range("A78").Value2 = "TOTAL"
$ws.Range("H78").Value2 = 561.96

Write-Host "Workbook updated through row 78."
